$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextValue 'D2' '22.103.67'
Set-TextValue 'D3' '1.550.98'
Set-TextValue 'E3' '  -1.32%  '
Set-TextValue 'D4' '1.000'
Set-TextValue 'E4' '  -0.16%  '
Set-TextValue 'D5' '0.9998'
Set-TextValue 'E5' '  -0.15%  '
Set-TextValue 'D6' '287.06'
Set-TextValue 'E6' '  -0.41%  '
Set-TextValue 'D7' '0.3805'
Set-TextValue 'E7' '  +2.36%  '
Set-TextValue 'D8' '0.3275'
Set-TextValue 'E8' '  -1.25%  '
Set-TextValue 'D9' '43.30'
Set-TextValue 'E9' '  -10.40%  '
Set-TextValue 'D10' '1.128'
Set-TextValue 'E10' '  -0.50%  '
Set-TextValue 'D11' '0.07325'
Set-TextValue 'E11' '  -2.15%  '
Set-TextValue 'D12' '1.000'
Set-TextValue 'E12' '  -0.17%  '
Set-TextValue 'D13' '19.99'
Set-TextValue 'E13' '  -3.65%  '
Set-TextValue 'D14' '5.794'
Set-TextValue 'E14' '  -2.38%  '
Set-TextValue 'D15' '6.762'
Set-TextValue 'E15' '  -2.05%  '
Set-TextValue 'D16' '1.549.91'
Set-TextValue 'E16' '  -1.44%  '
Set-TextValue 'D17' '0.00001087'
Set-TextValue 'E17' '  -2.84%  '
Set-TextValue 'D18' '0.06589'
Set-TextValue 'E18' '  -2.25%  '
Set-TextValue 'D19' '85.52'
Set-TextValue 'E19' '  -2.60%  '
Set-TextValue 'D20' '0.9997'
Set-TextValue 'E20' '  -0.15%  '
Set-TextValue 'D21' '6.366'
Set-TextValue 'E21' '  +0.16%  '
Set-TextValue 'D22' '16.04'
Set-TextValue 'E22' '  -3.13%  '
Set-TextValue 'E23' '  -3.51%  '
Set-TextValue 'D24' '22.108.57'
Set-TextValue 'E24' '  -1.51%  '
Set-TextValue 'D25' '2.292'
Set-TextValue 'E25' '  -3.78%  '
Set-TextValue 'D26' '2.518'
Set-TextValue 'E26' '  -2.38%  '
Set-TextValue 'D27' '149.45'
Set-TextValue 'E27' '  -2.50%  '
Set-TextValue 'D28' '19.03'
Set-TextValue 'E28' '  -3.40%  '
Set-TextValue 'D29' '4.902'
Set-TextValue 'E29' '  -2.34%  '
Set-TextValue 'D30' '121.03'
Set-TextValue 'E30' '  -2.76%  '
Set-TextValue 'D31' '1.722.32'
Set-TextValue 'E31' '  -1.47%  '
Set-TextValue 'D32' '1.063'
Set-TextValue 'E32' '  +0.98%  '
Set-TextValue 'D33' '5.872'
Set-TextValue 'E33' '  -4.33%  '
Set-TextValue 'D34' '1.863'
Set-TextValue 'E34' '  -7.37%  '
Set-TextValue 'D35' '0.08208'
Set-TextValue 'E35' '  -1.38%  '
Set-TextValue 'D36' '9.249'
Set-TextValue 'E36' '  -5.46%  '
Set-TextValue 'D37' '0.02308'
Set-TextValue 'E37' '  -6.33%  '
Set-TextValue 'E38' '  -3.02%  '
Set-TextValue 'D39' '5.253'
Set-TextValue 'E39' '  -1.81%  '
Set-TextValue 'D40' '0.2148'
Set-TextValue 'E40' '  -5.40%  '
Set-TextValue 'D41' '1.243'
Set-TextValue 'E41' '  -3.74%  '
Set-TextValue 'E42' '  -2.88%  '
Set-TextValue 'D43' '0.9994'
Set-TextValue 'E43' '  -0.16%  '
Set-TextValue 'D44' '0.6009'
Set-TextValue 'D45' '13.58'
Set-TextValue 'E45' '  -2.37%  '
Set-TextValue 'D46' '3.727'
Set-TextValue 'E46' '  -1.23%  '
Set-TextValue 'D47' '0.5802'
Set-TextValue 'E47' '  -5.75%  '
Set-TextValue 'D48' '1.975'
Set-TextValue 'E48' '  -4.01%  '
Set-TextValue 'D49' '121.72'
Set-TextValue 'E49' '  -3.22%  '
Set-TextValue 'D50' '1.171'
Set-TextValue 'E50' '  -3.23%  '
Set-TextValue 'D51' '0.07007'
Set-TextValue 'E51' '  -2.91%  '
